# Insert a new header row above the existing data (pushes all data down one row)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1:1").Insert()

# Populate the new header row
$ws.Range("A1").Value = "Animal"
$ws.Range("B1").Value = "year 2011"
$ws.Range("C1").Value = "year 2012"
$ws.Range("D1").Value = "year 2013"
$ws.Range("E1").Value = "year 2014"
$ws.Range("F1").Value = "year 2015"

# Match the final selection recorded in the workbook
$ws.Range("K11").Select()
